$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 54 (shifts existing rows 54..181 down to 55..182)
$ws.Rows("54").Insert()

# Populate the newly inserted row 54 with the new record
$ws.Cells.Item(54, 1).Value2  = 3
$ws.Cells.Item(54, 2).Value2  = "Femacal de La Calera"
$ws.Cells.Item(54, 3).Value2  = "Coquimbo"
$ws.Cells.Item(54, 4).Value2  = 44581
$ws.Cells.Item(54, 5).Value2  = 5
$ws.Cells.Item(54, 6).Value2  = 100112010
$ws.Cells.Item(54, 7).Value2  = "Achicoria"
$ws.Cells.Item(54, 8).Value2  = "Sin especificar"
$ws.Cells.Item(54, 9).Value2  = "Primera"
$ws.Cells.Item(54, 10).Value2 = 130
$ws.Cells.Item(54, 11).Value2 = 5500
$ws.Cells.Item(54, 12).Value2 = 6000
$ws.Cells.Item(54, 13).Value2 = 5769
$ws.Cells.Item(54, 14).Value2 = "$/caja 16 unidades"
$ws.Cells.Item(54, 15).Value2 = "Provincia de Quillota"
$ws.Cells.Item(54, 16).Value2 = 361
$ws.Cells.Item(54, 17).Value2 = 16
$ws.Cells.Item(54, 18).Value2 = "Hortaliza"

# Make sure the date cell keeps the date number format used by the rest of column D
$ws.Cells.Item(54, 4).NumberFormat = $ws.Cells.Item(55, 4).NumberFormat
